$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 108.09091
$ws.Range("I9").Value = 120.125
$ws.Range("J9").Value = 76
$ws.Range("K9").Value = 120.125
$ws.Range("L9").Value = 76
$ws.Range("M9").Value = 48.875
$ws.Range("N9").Value = -414
$ws.Range("H33").Value = 115.77778
$ws.Range("I33").Value = 94.35294
$ws.Range("K33").Value = 94.35294
$ws.Range("M33").Value = 134.64706
$ws.Range("H88").Value = 1123306.4
$ws.Range("I88").Value = 599.5
$ws.Range("J88").Value = 1764853.1
$ws.Range("K88").Value = 599.5
$ws.Range("L88").Value = 1764853.1
$ws.Range("M88").Value = -193.5
$ws.Range("N88").Value = -1765665.1
$ws.Range("H91").Value = 1123306.4
$ws.Range("I91").Value = 599.5
$ws.Range("J91").Value = 1764853.1
$ws.Range("K91").Value = 599.5
$ws.Range("L91").Value = 1764853.1
$ws.Range("M91").Value = 804.5
$ws.Range("N91").Value = -1767661.1
$ws.Range("H116").Value = 3315.375
$ws.Range("I116").Value = 3362.5
$ws.Range("J116").Value = 3249.4
$ws.Range("K116").Value = 3362.5
$ws.Range("L116").Value = 3249.4
$ws.Range("M116").Value = 79.5
$ws.Range("N116").Value = -10133.4
$ws.Range("H132").Value = 8136295.5
$ws.Range("I132").Value = 9264042
$ws.Range("K132").Value = 27792126
$ws.Range("M132").Value = -27789596
$ws.Range("H138").Value = 1339.56
$ws.Range("J138").Value = 1523.75
$ws.Range("L138").Value = 4571.25
$ws.Range("N138").Value = -14851.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3494.5356
$ws.Range("I32").Value = 3311.6
$ws.Range("K32").Value = 3311.6
$ws.Range("M32").Value = -3024.6
$ws.Range("H45").Value = 1371.7142
$ws.Range("I45").Value = 1501.8182
$ws.Range("J45").Value = 894.6667
$ws.Range("K45").Value = 1501.8182
$ws.Range("L45").Value = 894.6667
$ws.Range("M45").Value = -1124.8182
$ws.Range("N45").Value = -1648.6667
$ws.Range("H122").Value = 1716.1428
$ws.Range("I122").Value = 1899.8
$ws.Range("J122").Value = 1257
$ws.Range("K122").Value = 5699.4
$ws.Range("L122").Value = 3771
$ws.Range("M122").Value = -3249.4
$ws.Range("N122").Value = -8671
$ws.Range("H132").Value = 2364.2173
$ws.Range("I132").Value = 1975.8334
$ws.Range("J132").Value = 3762.4
$ws.Range("K132").Value = 5927.5002
$ws.Range("L132").Value = 11287.2
$ws.Range("M132").Value = -3397.5002
$ws.Range("N132").Value = -16347.2

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1341.5883
$ws.Range("I31").Value = 1120.6
$ws.Range("J31").Value = 2999
$ws.Range("K31").Value = 1120.6
$ws.Range("L31").Value = 2999
$ws.Range("M31").Value = -825.5999999999999
$ws.Range("N31").Value = -3589
$ws.Range("H34").Value = 1341.5883
$ws.Range("I34").Value = 1120.6
$ws.Range("J34").Value = 2999
$ws.Range("K34").Value = 1120.6
$ws.Range("L34").Value = 2999
$ws.Range("M34").Value = -918.5999999999999
$ws.Range("N34").Value = -3403
$ws.Range("H57").Value = 17000
$ws.Range("J57").Value = 17000
$ws.Range("L57").Value = 17000
$ws.Range("N57").Value = -18120
$ws.Range("H58").Value = 1101.0286
$ws.Range("I58").Value = 845.2593000000001
$ws.Range("J58").Value = 1964.25
$ws.Range("K58").Value = 845.2593000000001
$ws.Range("L58").Value = 1964.25
$ws.Range("M58").Value = -642.2593000000001
$ws.Range("N58").Value = -2370.25
$ws.Range("H132").Value = 4524.7295
$ws.Range("I132").Value = 4970.154
$ws.Range("K132").Value = 14910.462
$ws.Range("M132").Value = -12380.462
$ws.Range("H136").Value = 1101.0286
$ws.Range("I136").Value = 845.2593000000001
$ws.Range("J136").Value = 1964.25
$ws.Range("K136").Value = 2535.7779
$ws.Range("L136").Value = 5892.75
$ws.Range("M136").Value = 14.22209999999995
$ws.Range("N136").Value = -10992.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19232098
$ws.Range("I131").Value = 125000500
$ws.Range("K131").Value = 375001500
$ws.Range("M131").Value = -374996460
$ws.Range("H137").Value = 6780.9287
$ws.Range("I137").Value = 2223.3333
$ws.Range("J137").Value = 8939.789000000001
$ws.Range("K137").Value = 6669.999899999999
$ws.Range("L137").Value = 26819.367
$ws.Range("M137").Value = -1569.999899999999
$ws.Range("N137").Value = -37019.367

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2539.5454
$ws.Range("J126").Value = 2283.3333
$ws.Range("L126").Value = 6849.999899999999
$ws.Range("N126").Value = -11789.9999
$ws.Range("H132").Value = 2900.0715
$ws.Range("I132").Value = 2400.3333
$ws.Range("K132").Value = 7200.999899999999
$ws.Range("M132").Value = -4670.999899999999
$ws.Range("H133").Value = 42499.668
$ws.Range("J133").Value = 42499.668
$ws.Range("L133").Value = 42499.668
$ws.Range("N133").Value = -52619.668

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1740.5
$ws.Range("I22").Value = 922.375
$ws.Range("J22").Value = 2831.3333
$ws.Range("K22").Value = 922.375
$ws.Range("L22").Value = 2831.3333
$ws.Range("M22").Value = -627.375
$ws.Range("N22").Value = -3421.3333
$ws.Range("H27").Value = 1740.5
$ws.Range("I27").Value = 922.375
$ws.Range("J27").Value = 2831.3333
$ws.Range("K27").Value = 922.375
$ws.Range("L27").Value = 2831.3333
$ws.Range("M27").Value = -815.375
$ws.Range("N27").Value = -3045.3333
$ws.Range("H68").Value = 1851.2
$ws.Range("I68").Value = 1651.75
$ws.Range("J68").Value = 2649
$ws.Range("K68").Value = 1651.75
$ws.Range("L68").Value = 2649
$ws.Range("M68").Value = -902.75
$ws.Range("N68").Value = -4147
$ws.Range("H71").Value = 1851.2
$ws.Range("I71").Value = 1651.75
$ws.Range("J71").Value = 2649
$ws.Range("K71").Value = 8258.75
$ws.Range("L71").Value = 13245
$ws.Range("M71").Value = -4514.75
$ws.Range("N71").Value = -20733
$ws.Range("H122").Value = 62501628
$ws.Range("I122").Value = 125000750
$ws.Range("K122").Value = 375002250
$ws.Range("M122").Value = -374999800
$ws.Range("H136").Value = 6489.9546
$ws.Range("J136").Value = 2450.2
$ws.Range("L136").Value = 7350.599999999999
$ws.Range("N136").Value = -12450.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 12000
$ws.Range("J70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 12000
$ws.Range("J73").Value = 12000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -14184
$ws.Range("H107").Value = 323.33334
$ws.Range("I107").Value = 398
$ws.Range("J107").Value = 211.33333
$ws.Range("K107").Value = 1194
$ws.Range("L107").Value = 633.99999
$ws.Range("M107").Value = 726
$ws.Range("N107").Value = -4473.99999
$ws.Range("H132").Value = 3045.5151
$ws.Range("I132").Value = 2895.2778
$ws.Range("J132").Value = 3225.8
$ws.Range("K132").Value = 8685.8334
$ws.Range("L132").Value = 9677.400000000001
$ws.Range("M132").Value = -6155.8334
$ws.Range("N132").Value = -14737.4
$ws.Range("H136").Value = 729.3333
$ws.Range("I136").Value = 722.9091
$ws.Range("J136").Value = 800
$ws.Range("K136").Value = 2168.7273
$ws.Range("L136").Value = 2400
$ws.Range("M136").Value = 381.2727
$ws.Range("N136").Value = -7500
